# Rebuild the taxa table: replace the two "Unidentified ... eukaryots" rows
# with a fuller "incertae sedis" taxonomy (Prokaryota kingdom, plus
# Flagellates / Unicells / Eukarotic picoplankton branches down to Species),
# each block separated by a blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Scientific name (A), Rank (C), Parent name (D)
$rows = @(
    @(2,  "Eukaryota",                                       "Kingdom", ""),
    @(3,  "Prokaryota",                                      "Kingdom", ""),

    @(5,  "Flagellates phylum incertae sedis",                "Phylum",  "Eukaryota"),
    @(6,  "Flagellates classes incertae sedis",                "Class",   "Flagellates phylum incertae sedis"),
    @(7,  "Flagellates ordines incertae sedis",                "Order",   "Flagellates classes incertae sedis"),
    @(8,  "Flagellates families incertae sedis",               "Family",  "Flagellates ordines incertae sedis"),
    @(9,  "Flagellates genera incertae sedis",                 "Genus",   "Flagellates families incertae sedis"),
    @(10, "Flagellates species incertae sedis",                "Species", "Flagellates genera incertae sedis"),

    @(12, "Unicells kingdom incertae sedis",                   "Kingdom", ""),
    @(13, "Unicells phylum incertae sedis",                    "Phylum",  "Unicells kingdom incertae sedis"),
    @(14, "Unicells classes incertae sedis",                   "Class",   "Unicells phylum incertae sedis"),
    @(15, "Unicells ordines incertae sedis",                   "Order",   "Unicells classes incertae sedis"),
    @(16, "Unicells families incertae sedis",                  "Family",  "Unicells ordines incertae sedis"),
    @(17, "Unicells genera incertae sedis",                    "Genus",   "Unicells families incertae sedis"),
    @(18, "Unicells species incertae sedis",                   "Species", "Unicells genera incertae sedis"),

    @(20, "Eukarotic picoplankton phylum incertae sedis",      "Phylum",  "Eukaryota"),
    @(21, "Eukarotic picoplankton classes incertae sedis",     "Class",   "Eukarotic picoplankton phylum incertae sedis"),
    @(22, "Eukarotic picoplankton ordines incertae sedis",     "Order",   "Eukarotic picoplankton classes incertae sedis"),
    @(23, "Eukarotic picoplankton families incertae sedis",    "Family",  "Eukarotic picoplankton ordines incertae sedis"),
    @(24, "Eukarotic picoplankton genera incertae sedis",      "Genus",   "Eukarotic picoplankton families incertae sedis"),
    @(25, "Eukarotic picoplankton species incertae sedis",     "Species", "Eukarotic picoplankton genera incertae sedis")
)

# Clear out the old rows 3 and 4 first (they previously held the two
# "Unidentified ... eukaryots" entries, now replaced by the table above,
# which leaves some rows blank as separators between taxonomy blocks).
$ws.Range("A3:D4").Value = ""

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}
